$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 310 (D, E, F change; C and G stay the same) ---
$ws.Range("D310").Value = 107.3
$ws.Range("E310").Value = 105.1183
$ws.Range("F310").Value = 106.4

# --- Copy the date style from A310 down to A311:A313 so the new date cells ---
# --- keep the same number format (YYYY-MM-DD HH:MM:SS) as the rest of column A ---
$ws.Range("A310").Copy()
$ws.Range("A311:A313").PasteSpecial(-4122)

# --- New row 311 ---
$ws.Range("A311").Value = 45047.33333333334
$ws.Range("B311").Value = "FX_IDC:USDBDT"
$ws.Range("C311").Value = 106
$ws.Range("D311").Value = 108.2
$ws.Range("E311").Value = 104.72
$ws.Range("F311").Value = 107.4
$ws.Range("G311").Value = 0

# --- New row 312 ---
$ws.Range("A312").Value = 45078.33333333334
$ws.Range("B312").Value = "FX_IDC:USDBDT"
$ws.Range("C312").Value = 107.4
$ws.Range("D312").Value = 109.2
$ws.Range("E312").Value = 106.62
$ws.Range("F312").Value = 108.45
$ws.Range("G312").Value = 0

# --- New row 313 ---
$ws.Range("A313").Value = 45110.33333333334
$ws.Range("B313").Value = "FX_IDC:USDBDT"
$ws.Range("C313").Value = 108.45
$ws.Range("D313").Value = 108.6
$ws.Range("E313").Value = 107.38
$ws.Range("F313").Value = 108.55
$ws.Range("G313").Value = 0
